# Updates crypto price/volume data per the Wed Jan 3 02:33:15 UTC 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.465.33'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").Value = '2.376.99'
$ws.Range("E3").Value = '  -0.18%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = "'313.61"
$ws.Range("E5").Value = '  -1.36%  '

$ws.Range("D6").Value = "'108.27"
$ws.Range("E6").Value = '  -3.07%  '

$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = '  -1.19%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Value = "'0.612"
$ws.Range("E9").Value = '  -2.98%  '

$ws.Range("D10").Value = "'40.69"
$ws.Range("E10").Value = '  -3.83%  '

$ws.Range("D11").Value = "'0.0918"
$ws.Range("E11").Value = '  -1.07%  '

$ws.Range("D12").Value = "'8.51"
$ws.Range("E12").Value = '  -2.05%  '

$ws.Range("E13").Value = '  +0.58%  '

$ws.Range("D14").Value = "'0.982"
$ws.Range("E14").Value = '  -3.26%  '

$ws.Range("D15").Value = '2.738.03'
$ws.Range("E15").Value = '  -0.09%  '

$ws.Range("D16").Value = "'15.33"
$ws.Range("E16").Value = '  -3.11%  '

$ws.Range("D17").Value = '2.380.67'
$ws.Range("E17").Value = '  +0.92%  '

$ws.Range("D18").Value = '45.471.34'
$ws.Range("E18").Value = '  +0.61%  '

$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").Value = "'13.90"
$ws.Range("E19").Value = '  +5.38%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = "'7.25"
$ws.Range("E20").Value = '  -5.35%  '

$ws.Range("E21").Value = '  -1.23%  '

$ws.Range("D22").Value = "'73.43"
$ws.Range("E22").Value = '  -2.20%  '

$ws.Range("D23").Value = "'3.53"
$ws.Range("E23").Value = '  -1.27%  '

$ws.Range("D24").Value = "'259.89"
$ws.Range("E24").Value = '  -3.46%  '

$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = '  +2.20%  '

$ws.Range("E26").Value = '  -0.24%  '

$ws.Range("E27").Value = '  -0.89%  '

$ws.Range("D28").Value = "'7.25"
$ws.Range("E28").Value = '  -5.49%  '

$ws.Range("D29").Value = "'2.29"
$ws.Range("E29").Value = '  -1.38%  '

$ws.Range("D30").Value = "'0.0984"
$ws.Range("E30").Value = '  +5.42%  '

$ws.Range("D31").Value = "'22.38"
$ws.Range("E31").Value = '  -2.23%  '

$ws.Range("D32").Value = "'37.28"
$ws.Range("E32").Value = '  -6.25%  '

$ws.Range("D33").Value = "'166.53"
$ws.Range("E33").Value = '  -1.46%  '

$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = "'0.119"
$ws.Range("E35").Value = '  +1.05%  '

$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").Value = "'0.130"
$ws.Range("E36").Value = '  -2.02%  '

$ws.Range("D37").Value = "'4.70"
$ws.Range("E37").Value = '  -2.59%  '

$ws.Range("D38").Value = "'1.92"
$ws.Range("E38").Value = '  +9.72%  '

$ws.Range("D39").Value = "'4.00"
$ws.Range("E39").Value = '  +2.77%  '

$ws.Range("D40").Value = "'2.97"
$ws.Range("E40").Value = '  -1.48%  '

$ws.Range("D41").Value = "'0.0356"
$ws.Range("E41").Value = '  -3.07%  '

$ws.Range("D42").Value = "'98.62"
$ws.Range("E42").Value = '  -6.95%  '

$ws.Range("D43").Value = "'69.77"
$ws.Range("E43").Value = '  -2.77%  '

$ws.Range("D44").Value = "'0.229"
$ws.Range("E44").Value = '  -5.50%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").Value = "'12.72"
$ws.Range("E46").Value = '  -8.48%  '

$ws.Range("D47").Value = '1.818.55'
$ws.Range("E47").Value = '  +9.39%  '

$ws.Range("D48").Value = "'84.06"
$ws.Range("E48").Value = '  +5.56%  '

$ws.Range("D49").Value = "'5.80"
$ws.Range("E49").Value = '  +2.24%  '

$ws.Range("D50").Value = "'9.31"
$ws.Range("E50").Value = '  +2.45%  '

$ws.Range("D51").Value = "'111.04"
$ws.Range("E51").Value = '  -6.87%  '
